$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-21 Wednesday", "2026-01-22 Thursday"),
    @("176÷3=", "837÷3="),
    @("192÷9=", "801÷9="),
    @("952÷8=", "488÷6="),
    @("410÷8=", "145÷5="),
    @("604÷6=", "758÷6="),
    @("896÷8=", "134÷3="),
    @("961÷5=", "283÷9="),
    @("260÷3=", "772÷3="),
    @("934÷6=", "142÷5="),
    @("299÷5=", "794÷5="),
    @("211÷5=", "149÷7="),
    @("872÷2=", "270÷9="),
    @("134÷2=", "135÷8="),
    @("705÷5=", "373÷6="),
    @("100÷2=", "636÷9="),
    @("568÷4=", "677÷8="),
    @("961÷7=", "761÷4="),
    @("797÷5=", "274÷9="),
    @("222÷3=", "294÷4="),
    @("126÷7=", "313÷4="),
    @("395÷4=", "837÷2="),
    @("525÷3=", "709÷9="),
    @("334÷4=", "508÷2="),
    @("376÷3=", "841÷8="),
    @("498÷4=", "228÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
